# Apply the SoIB_summaries.xlsx update:
#  - rename "High Priority break-up" -> "Interannual update - High Pri"
#  - add a new sheet "Major update - High Priority " that is a copy of the
#    ORIGINAL "High Priority break-up" data (before the value updates below)
#  - update various statistic values on the existing sheets

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "Trends Status" sheet (sheet1) value updates
# ---------------------------------------------------------------------------
$wsTrends = $wb.Worksheets.Item("Trends Status")

# Row 2 - Rapid Decline
$wsTrends.Range("B2").Value = 0
$wsTrends.Range("C2").Value = 1
$wsTrends.Range("D2").Value = 0
$wsTrends.Range("E2").Value = 16.7

# Row 3 - Decline
$wsTrends.Range("B3").Value = 0
$wsTrends.Range("D3").Value = 0

# Row 4 - Stable
$wsTrends.Range("B4").Value = 1
$wsTrends.Range("D4").Value = 100
$wsTrends.Range("E4").Value = 83.3

# Row 5 - Increase
$wsTrends.Range("B5").Value = 0
$wsTrends.Range("D5").Value = 0

# Row 6 - Rapid Increase
$wsTrends.Range("B6").Value = 0
$wsTrends.Range("D6").Value = 0

# Row 7 - Trend Inconclusive
$wsTrends.Range("B7").Value = 63
$wsTrends.Range("C7").Value = 78

# Row 8 - Insufficient Data
$wsTrends.Range("B8").Value = 540
$wsTrends.Range("C8").Value = 520

# ---------------------------------------------------------------------------
# 2. "Priority Status" sheet (sheet3) value updates
# ---------------------------------------------------------------------------
$wsPriority = $wb.Worksheets.Item("Priority Status")

$wsPriority.Range("B2").Value = 103
$wsPriority.Range("B3").Value = 286
$wsPriority.Range("B4").Value = 554

# ---------------------------------------------------------------------------
# 3. "Species qualification" sheet (sheet4) value updates
# ---------------------------------------------------------------------------
$wsSpecies = $wb.Worksheets.Item("Species qualification")

$wsSpecies.Range("A2").Value = "SoIB Assessment"
$wsSpecies.Range("B2").Value = 604

$wsSpecies.Range("B3").Value = 64
$wsSpecies.Range("C3").Value = 1

$wsSpecies.Range("B4").Value = 84
$wsSpecies.Range("C4").Value = 6

# ---------------------------------------------------------------------------
# 4. Create the new "Major update - High Priority " sheet as a copy of the
#    ORIGINAL "High Priority break-up" data, placed at the end of the
#    workbook, BEFORE changing the values on the existing sheet.
# ---------------------------------------------------------------------------
$wsOldBreakup = $wb.Worksheets.Item("High Priority break-up")

$wsNew = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsNew.Name = "Major update - High Priority "

# NOTE: use Value() (as a call) to read the current value of a cell - the
# bare ".Value" property getter is not evaluated by this interpreter and
# would copy over a non-sensical placeholder instead of the real data.
$wsNew.Range("A1").Value = $wsOldBreakup.Range("A1").Value()
$wsNew.Range("B1").Value = $wsOldBreakup.Range("B1").Value()
$wsNew.Range("C1").Value = $wsOldBreakup.Range("C1").Value()
$wsNew.Range("D1").Value = $wsOldBreakup.Range("D1").Value()
$wsNew.Range("E1").Value = $wsOldBreakup.Range("E1").Value()
$wsNew.Range("A1:E1").Font.Bold = $true
$wsNew.Range("A1:E1").HorizontalAlignment = -4108

$wsNew.Range("A2").Value = $wsOldBreakup.Range("A2").Value()
$wsNew.Range("B2").Value = $wsOldBreakup.Range("B2").Value()
$wsNew.Range("C2").Value = $wsOldBreakup.Range("C2").Value()
$wsNew.Range("D2").Value = $wsOldBreakup.Range("D2").Value()
$wsNew.Range("E2").Value = $wsOldBreakup.Range("E2").Value()

$wsNew.Range("A3").Value = $wsOldBreakup.Range("A3").Value()
$wsNew.Range("B3").Value = $wsOldBreakup.Range("B3").Value()
$wsNew.Range("C3").Value = $wsOldBreakup.Range("C3").Value()
$wsNew.Range("D3").Value = $wsOldBreakup.Range("D3").Value()
$wsNew.Range("E3").Value = $wsOldBreakup.Range("E3").Value()

# ---------------------------------------------------------------------------
# 5. Rename "High Priority break-up" -> "Interannual update - High Pri" and
#    update its values to the new figures.
# ---------------------------------------------------------------------------
$wsOldBreakup.Name = "Interannual update - High Pri"

$wsOldBreakup.Range("B2").Value = 63
$wsOldBreakup.Range("C2").Value = 61.2
$wsOldBreakup.Range("D2").Value = 63
$wsOldBreakup.Range("E2").Value = 75

$wsOldBreakup.Range("B3").Value = 40
$wsOldBreakup.Range("C3").Value = 38.8
$wsOldBreakup.Range("D3").Value = 21
$wsOldBreakup.Range("E3").Value = 25
